$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 represents the "single line" comment direction symbols.
# Columns B (2) through AF (32), and AH (34), AI (35) had their sign and
# cell style flipped from negative/"Good" to positive/"Neutral".
# Column AG (33, the "\n" newline column) is left untouched.

$row = 25
$cols = @()
$cols += 2..32
$cols += 34
$cols += 35

# Use an existing "Neutral"-styled cell (with the sheet's normal border /
# number format / alignment direct formatting) as the format source, so
# that the edited cells end up sharing the very same cell style (s="4")
# as the rest of the sheet instead of Excel fabricating a brand new,
# stripped-down style entry.
$styleSource = $ws.Cells.Item(24, 12)   # L24, already styled "Neutral"
$styleSource.Copy() | Out-Null

foreach ($col in $cols) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $cell.Value = 23
}

$excel.CutCopyMode = 0
